$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B4").Value = 62210
$ws.Range("B16").Value = 53750
$ws.Range("B43").Value = 49240
$ws.Range("B54").Value = 52170
$ws.Range("B57").Value = 54300

$ws.Range("B57").Select()
